$wb = $excel.ActiveWorkbook

# --- Sheet: Statistics (L14) ---
$ws1 = $wb.Worksheets.Item("Statistics (L14)")
$ws1.Range("R3").Value = 355.4285714285714
$ws1.Range("R4").Value = 114.1959950030003
$ws1.Range("R5").Value = 241
$ws1.Range("R6").Value = 258.75
$ws1.Range("R7").Value = 323
$ws1.Range("R8").Value = 417.5
$ws1.Range("R9").Value = 619
$ws1.Range("R10").Value = 13040.72527472528
$ws1.Range("R11").Value = 1.063948529860846
$ws1.Range("R12").Value = 0.5808099381866265
$ws1.Range("R13").Value = 89.26530612244899

# --- Sheet: Statistics (All) ---
$ws2 = $wb.Worksheets.Item("Statistics (All)")
$ws2.Range("R3").Value = 217.5315315315315
$ws2.Range("R4").Value = 171.7247702691627
$ws2.Range("R6").Value = 66
$ws2.Range("R7").Value = 210
$ws2.Range("R8").Value = 323
$ws2.Range("R9").Value = 619
$ws2.Range("R10").Value = 29489.39672399672
$ws2.Range("R11").Value = 0.4461600987695634
$ws2.Range("R12").Value = -0.766514674913092
$ws2.Range("R13").Value = 137.2688905121338

# --- Sheet: Kosovo Raw Data ---
$ws3 = $wb.Worksheets.Item("Kosovo Raw Data")
$ws3.Range("T27").Value = 32
$ws3.Range("T28").Value = 60
$ws3.Range("T29").Value = 62
$ws3.Range("T30").Value = 70
$ws3.Range("T31").Value = 84
$ws3.Range("T32").Value = 86
$ws3.Range("T33").Value = 89
$ws3.Range("T34").Value = 92
$ws3.Range("T35").Value = 104
$ws3.Range("T36").Value = 105
$ws3.Range("T37").Value = 114
$ws3.Range("T38").Value = 115
$ws3.Range("T39").Value = 118
$ws3.Range("T40").Value = 116
$ws3.Range("T41").Value = 121
$ws3.Range("T42").Value = 138
$ws3.Range("T43").Value = 149
$ws3.Range("T44").Value = 181
$ws3.Range("T45").Value = 182
$ws3.Range("T46").Value = 191
$ws3.Range("T47").Value = 218
$ws3.Range("T48").Value = 296
$ws3.Range("T49").Value = 306
$ws3.Range("T50").Value = 313
$ws3.Range("T51").Value = 343
$ws3.Range("T52").Value = 359
$ws3.Range("T53").Value = 384
$ws3.Range("T54").Value = 405
$ws3.Range("T55").Value = 447
$ws3.Range("T56").Value = 460
$ws3.Range("T57").Value = 461
$ws3.Range("T58").Value = 474
$ws3.Range("T59").Value = 491
$ws3.Range("T60").Value = 522
$ws3.Range("T61").Value = 546
$ws3.Range("T62").Value = 576
$ws3.Range("T63").Value = 557
$ws3.Range("T64").Value = 536
$ws3.Range("T65").Value = 528
$ws3.Range("T66").Value = 513
$ws3.Range("T67").Value = 493
$ws3.Range("T68").Value = 465
$ws3.Range("T69").Value = 445
$ws3.Range("T70").Value = 426
$ws3.Range("T71").Value = 340
$ws3.Range("T72").Value = 301
$ws3.Range("T73").Value = 273
$ws3.Range("T74").Value = 212
$ws3.Range("T75").Value = 189
$ws3.Range("T76").Value = 201
$ws3.Range("T77").Value = 210
$ws3.Range("T78").Value = 219
$ws3.Range("T79").Value = 227
$ws3.Range("T80").Value = 226
$ws3.Range("T81").Value = 235
$ws3.Range("T82").Value = 236
$ws3.Range("T83").Value = 220
$ws3.Range("T84").Value = 205
$ws3.Range("T85").Value = 191
$ws3.Range("T86").Value = 202
$ws3.Range("T87").Value = 203
$ws3.Range("T88").Value = 214
$ws3.Range("T89").Value = 218
$ws3.Range("T90").Value = 214
$ws3.Range("T91").Value = 217
$ws3.Range("T92").Value = 223
$ws3.Range("T93").Value = 217
$ws3.Range("T94").Value = 202
$ws3.Range("T95").Value = 205
$ws3.Range("T96").Value = 199
$ws3.Range("T97").Value = 210
$ws3.Range("T98").Value = 237
$ws3.Range("T99").Value = 246
$ws3.Range("T100").Value = 241
$ws3.Range("T101").Value = 243
$ws3.Range("T102").Value = 252
$ws3.Range("T103").Value = 279
$ws3.Range("T104").Value = 313
$ws3.Range("T105").Value = 320
$ws3.Range("T106").Value = 326
$ws3.Range("T107").Value = 354
$ws3.Range("T108").Value = 374
$ws3.Range("T109").Value = 432
$ws3.Range("T110").Value = 477
$ws3.Range("T111").Value = 500
$ws3.Range("T112").Value = 619
